$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.08059215988451403
$ws.Range("J2").Value = 0.08059215988451404
$ws.Range("M2").Value = 0.032838
$ws.Range("N2").Value = 0.098514
$ws.Range("O2").Value = 0.007146324094219707
$ws.Range("P2").Value = 0.007146324094219707
$ws.Range("Q2").Value = 0.006295504332000001
$ws.Range("R2").Value = 0.056659538988
$ws.Range("S2").Value = 0.0005759376939879095
$ws.Range("T2").Value = 0.0005759376939879097

# Row 3
$ws.Range("I3").Value = 0.08059215988451403
$ws.Range("J3").Value = 0.08059215988451404
$ws.Range("O3").Value = 0.03951718316124263
$ws.Range("P3").Value = 0.03951718316124263
$ws.Range("S3").Value = 0.003184775143516491
$ws.Range("T3").Value = 0.003184775143516492

# Row 4
$ws.Range("I4").Value = 0.08059215988451403
$ws.Range("J4").Value = 0.08059215988451404
$ws.Range("M4").Value = 3.814633
$ws.Range("N4").Value = 11.443899
$ws.Range("O4").Value = 0.8301542030119253
$ws.Range("P4").Value = 0.8301542030119253
$ws.Range("Q4").Value = 0.7313185509620002
$ws.Range("R4").Value = 6.581866958658001
$ws.Range("S4").Value = 0.0669039202579384
$ws.Range("T4").Value = 0.06690392025793841

# Row 5
$ws.Range("I5").Value = 0.08059215988451403
$ws.Range("J5").Value = 0.08059215988451404
$ws.Range("M5").Value = 0.5660336666666667
$ws.Range("N5").Value = 1.698101
$ws.Range("O5").Value = 0.1231822897326124
$ws.Range("P5").Value = 0.1231822897326124
$ws.Range("Q5").Value = 0.1085165783713333
$ws.Range("R5").Value = 0.976649205342
$ws.Range("S5").Value = 0.00992752678907123
$ws.Range("T5").Value = 0.00992752678907123

# Row 6
$ws.Range("G6").Value = 2.187103
$ws.Range("H6").Value = 6.561309
$ws.Range("I6").Value = 0.9194078401154859
$ws.Range("J6").Value = 0.919407840115486
$ws.Range("M6").Value = 0.032838
$ws.Range("N6").Value = 0.098514
$ws.Range("O6").Value = 0.007146324094219707
$ws.Range("P6").Value = 0.007146324094219707
$ws.Range("Q6").Value = 0.071820088314
$ws.Range("R6").Value = 0.646380794826
$ws.Range("S6").Value = 0.006570386400231797
$ws.Range("T6").Value = 0.006570386400231798

# Row 7
$ws.Range("G7").Value = 2.187103
$ws.Range("H7").Value = 6.561309
$ws.Range("I7").Value = 0.9194078401154859
$ws.Range("J7").Value = 0.919407840115486
$ws.Range("O7").Value = 0.03951718316124263
$ws.Range("P7").Value = 0.03951718316124263
$ws.Range("Q7").Value = 0.397145098255
$ws.Range("R7").Value = 3.574305884295
$ws.Range("S7").Value = 0.03633240801772614
$ws.Range("T7").Value = 0.03633240801772614

# Row 8
$ws.Range("G8").Value = 2.187103
$ws.Range("H8").Value = 6.561309
$ws.Range("I8").Value = 0.9194078401154859
$ws.Range("J8").Value = 0.919407840115486
$ws.Range("M8").Value = 3.814633
$ws.Range("N8").Value = 11.443899
$ws.Range("O8").Value = 0.8301542030119253
$ws.Range("P8").Value = 0.8301542030119253
$ws.Range("Q8").Value = 8.342995278199
$ws.Range("R8").Value = 75.086957503791
$ws.Range("S8").Value = 0.7632502827539869
$ws.Range("T8").Value = 0.763250282753987

# Row 9
$ws.Range("G9").Value = 2.187103
$ws.Range("H9").Value = 6.561309
$ws.Range("I9").Value = 0.9194078401154859
$ws.Range("J9").Value = 0.919407840115486
$ws.Range("M9").Value = 0.5660336666666667
$ws.Range("N9").Value = 1.698101
$ws.Range("O9").Value = 0.1231822897326124
$ws.Range("P9").Value = 0.1231822897326124
$ws.Range("Q9").Value = 1.237973930467667
$ws.Range("R9").Value = 11.141765374209
$ws.Range("S9").Value = 0.1132547629435412
$ws.Range("T9").Value = 0.1132547629435412
